$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "M2"
$ws.Range("B2").Value = "C1qa"
$ws.Range("C2").Value = "Cspg4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 191.5877736666667
$ws.Range("H2").Value = 574.763321
$ws.Range("I2").Value = 0.9786143588951871
$ws.Range("J2").Value = 0.9786143588951871
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.723583666666666
$ws.Range("N2").Value = 14.170751
$ws.Range("O2").Value = 0.1402914168557431
$ws.Range("P2").Value = 0.1402914168557431
$ws.Range("Q2").Value = 904.9808784248967
$ws.Range("R2").Value = 8144.827905824071
$ws.Range("S2").Value = 0.1372911949647805
$ws.Range("T2").Value = 0.1372911949647805

$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "C1qa"
$ws.Range("C3").Value = "Cspg4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 191.5877736666667
$ws.Range("H3").Value = 574.763321
$ws.Range("I3").Value = 0.9786143588951871
$ws.Range("J3").Value = 0.9786143588951871
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.354856
$ws.Range("N3").Value = 31.064568
$ws.Range("O3").Value = 0.3075413758051058
$ws.Range("P3").Value = 0.3075413758051058
$ws.Range("Q3").Value = 1983.863807678925
$ws.Range("R3").Value = 17854.77426911033
$ws.Range("S3").Value = 0.3009644063172575
$ws.Range("T3").Value = 0.3009644063172575

$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "C1qa"
$ws.Range("C4").Value = "Cspg4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 191.5877736666667
$ws.Range("H4").Value = 574.763321
$ws.Range("I4").Value = 0.9786143588951871
$ws.Range("J4").Value = 0.9786143588951871
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6477360000000001
$ws.Range("N4").Value = 1.943208
$ws.Range("O4").Value = 0.01923789385371424
$ws.Range("P4").Value = 0.01923789385371424
$ws.Range("Q4").Value = 124.098298163752
$ws.Range("R4").Value = 1116.884683473768
$ws.Range("S4").Value = 0.01882647916014622
$ws.Range("T4").Value = 0.01882647916014622

$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "C1qa"
$ws.Range("C5").Value = "Cspg4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 191.5877736666667
$ws.Range("H5").Value = 574.763321
$ws.Range("I5").Value = 0.9786143588951871
$ws.Range("J5").Value = 0.9786143588951871
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.94362233333333
$ws.Range("N5").Value = 53.830867
$ws.Range("O5").Value = 0.5329293134854368
$ws.Range("P5").Value = 0.5329293134854368
$ws.Range("Q5").Value = 3437.778654358812
$ws.Range("R5").Value = 30940.00788922931
$ws.Range("S5").Value = 0.521532278453003
$ws.Range("T5").Value = 0.521532278453003

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "C1qa"
$ws.Range("C6").Value = "Cspg4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 4.186764
$ws.Range("H6").Value = 12.560292
$ws.Range("I6").Value = 0.02138564110481286
$ws.Range("J6").Value = 0.02138564110481286
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.723583666666666
$ws.Range("N6").Value = 14.170751
$ws.Range("O6").Value = 0.1402914168557431
$ws.Range("P6").Value = 0.1402914168557431
$ws.Range("Q6").Value = 19.776530046588
$ws.Range("R6").Value = 177.988770419292
$ws.Range("S6").Value = 0.003000221890962616
$ws.Range("T6").Value = 0.003000221890962616

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "C1qa"
$ws.Range("C7").Value = "Cspg4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 4.186764
$ws.Range("H7").Value = 12.560292
$ws.Range("I7").Value = 0.02138564110481286
$ws.Range("J7").Value = 0.02138564110481286
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.354856
$ws.Range("N7").Value = 31.064568
$ws.Range("O7").Value = 0.3075413758051058
$ws.Range("P7").Value = 0.3075413758051058
$ws.Range("Q7").Value = 43.353338325984
$ws.Range("R7").Value = 390.180044933856
$ws.Range("S7").Value = 0.00657696948784837
$ws.Range("T7").Value = 0.00657696948784837

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "C1qa"
$ws.Range("C8").Value = "Cspg4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 4.186764
$ws.Range("H8").Value = 12.560292
$ws.Range("I8").Value = 0.02138564110481286
$ws.Range("J8").Value = 0.02138564110481286
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.6477360000000001
$ws.Range("N8").Value = 1.943208
$ws.Range("O8").Value = 0.01923789385371424
$ws.Range("P8").Value = 0.01923789385371424
$ws.Range("Q8").Value = 2.711917766304
$ws.Range("R8").Value = 24.407259896736
$ws.Range("S8").Value = 0.0004114146935680179
$ws.Range("T8").Value = 0.0004114146935680179

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "C1qa"
$ws.Range("C9").Value = "Cspg4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 4.186764
$ws.Range("H9").Value = 12.560292
$ws.Range("I9").Value = 0.02138564110481286
$ws.Range("J9").Value = 0.02138564110481286
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.94362233333333
$ws.Range("N9").Value = 53.830867
$ws.Range("O9").Value = 0.5329293134854368
$ws.Range("P9").Value = 0.5329293134854368
$ws.Range("Q9").Value = 75.12571201479601
$ws.Range("R9").Value = 676.131408133164
$ws.Range("S9").Value = 0.01139703503243386
$ws.Range("T9").Value = 0.01139703503243386
